# Apply the "Add files via upload" edit: append 9 new dictionary rows
# (rows 182-190) to Sheet1, matching the shared-strings / cell layout
# introduced in the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper source cells whose existing cell styles we reuse (via
# PasteSpecial-formats) so that no new/duplicate font or cellXf entries
# get created in styles.xml - this mirrors the pre-existing style
# palette already used throughout the sheet.
#   style 2  (B52)  -> Lucida Sans Unicode, 10pt, #333333        (IPA cells)
#   style 3  (G52)  -> default font + wrapText                   (2-line examples)
#   style 5  (E57)  -> Source Sans Pro, 10pt, #333333             (definitions)
#   style 6  (G57)  -> italic Inherit 10pt #333333 + wrap/indent  (collocations)
#   style 26 (A158) -> bold Calibri 11pt + wrap, vertical center  (headwords)
# ---------------------------------------------------------------------
$styleA = $ws.Range("A158")   # s=26 headword style
$styleB = $ws.Range("B52")    # s=2  IPA style
$styleE = $ws.Range("E57")    # s=5  definition style
$styleF6 = $ws.Range("G57")   # s=6  italic indent style
$styleF3 = $ws.Range("G52")   # s=3  wrap style

function Set-HeadwordStyle($cell) {
    $styleA.Copy()
    $cell.PasteSpecial(-4122)
}
function Set-IpaStyle($cell) {
    $styleB.Copy()
    $cell.PasteSpecial(-4122)
}
function Set-DefStyle($cell) {
    $styleE.Copy()
    $cell.PasteSpecial(-4122)
}
function Set-IndentItalicStyle($cell) {
    $styleF6.Copy()
    $cell.PasteSpecial(-4122)
}
function Set-WrapStyle($cell) {
    $styleF3.Copy()
    $cell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Row 182 - strain
# ---------------------------------------------------------------------
$ws.Range("A182").Value = "strain"
Set-HeadwordStyle $ws.Range("A182")

$ws.Range("E182").Value = "mental pressure or worry felt by somebody because they have too much to do or manage; something that causes this pressure"
Set-DefStyle $ws.Range("E182")

$ws.Range("B182").Value = "/streɪn/"
Set-IpaStyle $ws.Range("B182")

$ws.Range("F182").Value = "(under strain) Their marriage is under great strain at the moment.`n (strain on something) These repayments are putting a strain on our finances."
Set-WrapStyle $ws.Range("F182")

$ws.Range("G182").Value = "The transport service cannot cope with the strain of so many additional passengers."
Set-IndentItalicStyle $ws.Range("G182")

$ws.Range("C182").Value = "n"

$ws.Rows.Item(182).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 183 - alternative
# ---------------------------------------------------------------------
$ws.Range("A183").Value = "alternative"
Set-HeadwordStyle $ws.Range("A183")

$ws.Range("B183").Value = "/ɔːlˈtɜːnətɪv/"
Set-IpaStyle $ws.Range("B183")

$ws.Range("E183").Value = "a thing that you can choose to do or have out of two or more possibilities"
Set-DefStyle $ws.Range("E183")

$rF183 = $ws.Range("F183")
$rF183.Value = "to offer/provide an alternative"
Set-IndentItalicStyle $rF183
$rF183.Value = "to offer/provide an alternative"
$runF183 = $rF183.Characters(4, 28)
$runF183.Font.Italic = $true
$runF183.Font.Size = 10
$runF183.Font.Color = 3355443
$runF183.Font.Name = "Inherit"

$rG183 = $ws.Range("G183")
$rG183.Value = "(alternative for somebody) This treatment is the only alternative for some patients."
Set-DefStyle $rG183
$rG183.Value = "(alternative for somebody) This treatment is the only alternative for some patients."
$runG183a = $rG183.Characters(27, 23)
$runG183a.Font.Italic = $true
$runG183a.Font.Size = 10
$runG183a.Font.Color = 3355443
$runG183a.Font.Name = "Source Sans Pro"
$runG183b = $rG183.Characters(50, 16)
$runG183b.Font.Italic = $true
$runG183b.Font.Size = 10
$runG183b.Font.Color = 3355443
$runG183b.Font.Name = "Inherit"
$runG183c = $rG183.Characters(66, 19)
$runG183c.Font.Italic = $true
$runG183c.Font.Size = 10
$runG183c.Font.Color = 3355443
$runG183c.Font.Name = "Source Sans Pro"

$ws.Range("C183").Value = "n"

# ---------------------------------------------------------------------
# Row 184 - substitude
# ---------------------------------------------------------------------
$ws.Range("A184").Value = "substitude"
Set-HeadwordStyle $ws.Range("A184")

$ws.Range("B184").Value = "/ˈsʌbstɪtjuːt/"
Set-IpaStyle $ws.Range("B184")

$ws.Range("E184").Value = "a person or thing that you use or have instead of the one you normally use or have"
Set-DefStyle $ws.Range("E184")

$ws.Range("F184").Value = "substitute for somebody/something"
Set-DefStyle $ws.Range("F184")

$ws.Range("C184").Value = "n"
$ws.Range("D184").Value = "alternative"

# ---------------------------------------------------------------------
# Row 185 - levy a tax
# ---------------------------------------------------------------------
$ws.Range("A185").Value = "levy a tax"
Set-HeadwordStyle $ws.Range("A185")

$ws.Range("E185").Value = "đánh thuế"
Set-DefStyle $ws.Range("E185")

$ws.Range("B185").Value = "/ˈlevi/"
Set-IpaStyle $ws.Range("B185")

# ---------------------------------------------------------------------
# Row 186 - significant hurdle
# ---------------------------------------------------------------------
$ws.Range("A186").Value = "significant hurdle "

$ws.Range("E186").Value = "a problem or difficulty that must be solved or dealt with before you can achieve something"
Set-DefStyle $ws.Range("E186")

$ws.Range("D186").Value = "obstacle"

$ws.Range("B186").Value = "/ˈhɜːdl/"
Set-IpaStyle $ws.Range("B186")

$ws.Range("C186").Value = "n"

# ---------------------------------------------------------------------
# Row 187 - buzz about
# ---------------------------------------------------------------------
$ws.Range("A187").Value = "buzz about"
Set-HeadwordStyle $ws.Range("A187")

$ws.Range("E187").Value = "xôn xao về"
Set-DefStyle $ws.Range("E187")

$ws.Range("G187").Value = "the internet is buzzing about the tiny house of Elon Musk  "

# ---------------------------------------------------------------------
# Row 188 - furnishings
# ---------------------------------------------------------------------
$ws.Range("F188").Value = "bedroom furnishings"

$ws.Range("E188").Value = "the furniture, carpets, curtains, etc. in a room or house"
Set-DefStyle $ws.Range("E188")

$ws.Range("B188").Value = "/ˈfɜːrnɪʃɪŋz/"
Set-IpaStyle $ws.Range("B188")

$ws.Range("A188").Value = "furnishings"

$ws.Range("C188").Value = "n"

# ---------------------------------------------------------------------
# Row 189 - collobrate
# ---------------------------------------------------------------------
$ws.Range("A189").Value = "collobrate"
Set-HeadwordStyle $ws.Range("A189")

$ws.Range("D189").Value = "cooperate"

# ---------------------------------------------------------------------
# Row 190 - expertise
# ---------------------------------------------------------------------
$ws.Range("A190").Value = "expertise (n) <> expert (n/a)"

$ws.Range("E190").Value = "có chuyên môn <> chuyên gia"
Set-DefStyle $ws.Range("E190")

$ws.Range("F190").Value = "expert opinion"

$ws.Range("G190").Value = "he has expertise in artificial inteligent"

$ws.Application.CutCopyMode = $false
